$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# F1: estadoUnidad -> estadoEntrega
$ws.Range("F1").Value = "estadoEntrega"
# G1: estadoEntrega -> comment (column H's old header moves into G)
$ws.Range("G1").Value = "comment"
# H1: comment -> removed (column H is dropped entirely)
$ws.Range("H1").ClearContents()

# --- Data row (row 2) ---
$ws.Range("A2").Value = "RIKO VIX "

# B2 must stay a text value ("123"), not be auto-converted to a number
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "123"
$ws.Range("B2").ClearFormats()

$ws.Range("C2").Value = "CDMX"
$ws.Range("D2").Value = "Bodega"
$ws.Range("E2").Value = "Correctivo"
$ws.Range("F2").Value = "Estaba por entregar"
$ws.Range("G2").Value = "se chingo una llanta"
# H2: comment value removed (column H is dropped entirely)
$ws.Range("H2").ClearContents()
